# Auto-generated edit script applying numeric corrections per commit diff
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC (index 1) ----
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(12, 8).Value = 171.75  # H12: 170.5 -> 171.75
$ws.Cells.Item(12, 9).Value = 171.75  # I12: 170.5 -> 171.75
$ws.Cells.Item(12, 11).Value = 171.75  # K12: 170.5 -> 171.75
$ws.Cells.Item(12, 13).Value = -1.75  # M12: -0.5 -> -1.75
$ws.Cells.Item(33, 8).Value = 330.42856  # H33: 406.9091 -> 330.42856
$ws.Cells.Item(33, 9).Value = 330.42856  # I33: 406.9091 -> 330.42856
$ws.Cells.Item(33, 11).Value = 330.42856  # K33: 406.9091 -> 330.42856
$ws.Cells.Item(33, 13).Value = -101.42856  # M33: -177.9091 -> -101.42856
$ws.Cells.Item(43, 8).Value = 2448.8  # H43: 1789 -> 2448.8
$ws.Cells.Item(43, 10).Value = 2883  # J43: 1784.75 -> 2883
$ws.Cells.Item(43, 12).Value = 2883  # L43: 1784.75 -> 2883
$ws.Cells.Item(43, 14).Value = -3021  # N43: -1922.75 -> -3021
$ws.Cells.Item(64, 8).Value = 3983.5483  # H64: 3968.4375 -> 3983.5483
$ws.Cells.Item(64, 10).Value = 5276.6665  # J64: 5099 -> 5276.6665
$ws.Cells.Item(64, 12).Value = 5276.6665  # L64: 5099 -> 5276.6665
$ws.Cells.Item(64, 14).Value = -5772.6665  # N64: -5595 -> -5772.6665
$ws.Cells.Item(67, 8).Value = 3983.5483  # H67: 3968.4375 -> 3983.5483
$ws.Cells.Item(67, 10).Value = 5276.6665  # J67: 5099 -> 5276.6665
$ws.Cells.Item(67, 12).Value = 5276.6665  # L67: 5099 -> 5276.6665
$ws.Cells.Item(67, 14).Value = -6992.6665  # N67: -6815 -> -6992.6665
$ws.Cells.Item(69, 8).Value = 11514.667  # H69: 13686.556 -> 11514.667
$ws.Cells.Item(69, 9).Value = 11882.714  # I69: 13030 -> 11882.714
$ws.Cells.Item(69, 10).Value = 10999.4  # J69: 14999.667 -> 10999.4
$ws.Cells.Item(69, 11).Value = 35648.142  # K69: 39090 -> 35648.142
$ws.Cells.Item(69, 12).Value = 32998.2  # L69: 44999.001 -> 32998.2
$ws.Cells.Item(69, 13).Value = -34774.142  # M69: -38216 -> -34774.142
$ws.Cells.Item(69, 14).Value = -34746.2  # N69: -46747.001 -> -34746.2
$ws.Cells.Item(72, 8).Value = 11514.667  # H72: 13686.556 -> 11514.667
$ws.Cells.Item(72, 9).Value = 11882.714  # I72: 13030 -> 11882.714
$ws.Cells.Item(72, 10).Value = 10999.4  # J72: 14999.667 -> 10999.4
$ws.Cells.Item(72, 11).Value = 106944.426  # K72: 117270 -> 106944.426
$ws.Cells.Item(72, 12).Value = 98994.59999999999  # L72: 134997.003 -> 98994.59999999999
$ws.Cells.Item(72, 13).Value = -102576.426  # M72: -112902 -> -102576.426
$ws.Cells.Item(72, 14).Value = -107730.6  # N72: -143733.003 -> -107730.6
$ws.Cells.Item(74, 8).Value = 12083.333  # H74: 8859.091 -> 12083.333
$ws.Cells.Item(74, 9).Value = 8500  # I74: 6745 -> 8500
$ws.Cells.Item(74, 11).Value = 8500  # K74: 6745 -> 8500
$ws.Cells.Item(74, 13).Value = -7564  # M74: -5809 -> -7564
$ws.Cells.Item(77, 8).Value = 12083.333  # H77: 8859.091 -> 12083.333
$ws.Cells.Item(77, 9).Value = 8500  # I77: 6745 -> 8500
$ws.Cells.Item(77, 11).Value = 42500  # K77: 33725 -> 42500
$ws.Cells.Item(77, 13).Value = -37820  # M77: -29045 -> -37820
$ws.Cells.Item(86, 8).Value = 12083.167  # H86: 13600 -> 12083.167
$ws.Cells.Item(86, 9).Value = 4499.8  # I86: 4500 -> 4499.8
$ws.Cells.Item(86, 11).Value = 4499.8  # K86: 4500 -> 4499.8
$ws.Cells.Item(86, 13).Value = -3376.8  # M86: -3377 -> -3376.8
$ws.Cells.Item(89, 8).Value = 12083.167  # H89: 13600 -> 12083.167
$ws.Cells.Item(89, 9).Value = 4499.8  # I89: 4500 -> 4499.8
$ws.Cells.Item(89, 11).Value = 22499  # K89: 22500 -> 22499
$ws.Cells.Item(89, 13).Value = -16883  # M89: -16884 -> -16883
$ws.Cells.Item(116, 8).Value = 16836  # H116: 17391.666 -> 16836
$ws.Cells.Item(116, 9).Value = 11503.4  # I116: 12392.667 -> 11503.4
$ws.Cells.Item(116, 11).Value = 11503.4  # K116: 12392.667 -> 11503.4
$ws.Cells.Item(116, 13).Value = -8061.4  # M116: -8950.666999999999 -> -8061.4
$ws.Cells.Item(131, 8).Value = 2629.8  # H131: 2691.1667 -> 2629.8
$ws.Cells.Item(131, 9).Value = 2629.8  # I131: 2691.1667 -> 2629.8
$ws.Cells.Item(131, 11).Value = 7889.400000000001  # K131: 8073.500100000001 -> 7889.400000000001
$ws.Cells.Item(131, 13).Value = -2849.400000000001  # M131: -3033.500100000001 -> -2849.400000000001
$ws.Cells.Item(132, 8).Value = 18479.924  # H132: 32475.572 -> 18479.924
$ws.Cells.Item(132, 9).Value = 21249.727  # I132: 44167.4 -> 21249.727
$ws.Cells.Item(132, 11).Value = 63749.181  # K132: 132502.2 -> 63749.181
$ws.Cells.Item(132, 13).Value = -61219.181  # M132: -129972.2 -> -61219.181
$ws.Cells.Item(135, 8).Value = 3001.8  # H135: 2179.4285 -> 3001.8
$ws.Cells.Item(135, 9).Value = 2842.7778  # I135: 2006.0769 -> 2842.7778
$ws.Cells.Item(135, 11).Value = 25585.0002  # K135: 18054.6921 -> 25585.0002
$ws.Cells.Item(135, 13).Value = -23050.0002  # M135: -15519.6921 -> -23050.0002
$ws.Cells.Item(137, 8).Value = 3762.4546  # H137: 3639.9565 -> 3762.4546
$ws.Cells.Item(137, 9).Value = 1542.6765  # I137: 1509.4722 -> 1542.6765
$ws.Cells.Item(137, 11).Value = 4628.029500000001  # K137: 4528.4166 -> 4628.029500000001
$ws.Cells.Item(137, 13).Value = -2078.029500000001  # M137: -1978.4166 -> -2078.029500000001
$ws.Cells.Item(138, 8).Value = 4955.66  # H138: 5014.745 -> 4955.66
$ws.Cells.Item(138, 10).Value = 4692.9375  # J138: 4775.8667 -> 4692.9375
$ws.Cells.Item(138, 12).Value = 14078.8125  # L138: 14327.6001 -> 14078.8125
$ws.Cells.Item(138, 14).Value = -24358.8125  # N138: -24607.6001 -> -24358.8125

# ---- Sheet: ARM (index 2) ----
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(32, 8).Value = 2438.36  # H32: 3437.68 -> 2438.36
$ws.Cells.Item(32, 9).Value = 2291.0957  # I32: 2315 -> 2291.0957
$ws.Cells.Item(32, 10).Value = 4745.5  # J32: 18353.285 -> 4745.5
$ws.Cells.Item(32, 11).Value = 2291.0957  # K32: 2315 -> 2291.0957
$ws.Cells.Item(32, 12).Value = 4745.5  # L32: 18353.285 -> 4745.5
$ws.Cells.Item(32, 13).Value = -2004.0957  # M32: -2028 -> -2004.0957
$ws.Cells.Item(32, 14).Value = -5319.5  # N32: -18927.285 -> -5319.5
$ws.Cells.Item(40, 8).Value = 0  # H40: 20000 -> 0
$ws.Cells.Item(40, 10).Value = 0  # J40: 20000 -> 0
$ws.Cells.Item(40, 12).ClearContents()  # L40: delete (was 20000)
$ws.Cells.Item(40, 14).Value = 0  # N40: -20352 -> 0
$ws.Cells.Item(45, 8).Value = 4071.1428  # H45: 4071.2856 -> 4071.1428
$ws.Cells.Item(45, 9).Value = 2624.75  # I45: 2625 -> 2624.75
$ws.Cells.Item(45, 11).Value = 2624.75  # K45: 2625 -> 2624.75
$ws.Cells.Item(45, 13).Value = -2247.75  # M45: -2248 -> -2247.75
$ws.Cells.Item(132, 8).Value = 3047.3455  # H132: 3072.818 -> 3047.3455
$ws.Cells.Item(132, 9).Value = 2205  # I132: 2238.3572 -> 2205
$ws.Cells.Item(132, 11).Value = 6615  # K132: 6715.071599999999 -> 6615
$ws.Cells.Item(132, 13).Value = -4085  # M132: -4185.071599999999 -> -4085

# ---- Sheet: BSM (index 3) ----
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(20, 8).Value = 49257.137  # H20: 54062.85 -> 49257.137
$ws.Cells.Item(20, 10).Value = 2761.125  # J20: 3281.5 -> 2761.125
$ws.Cells.Item(20, 12).Value = 2761.125  # L20: 3281.5 -> 2761.125
$ws.Cells.Item(20, 14).Value = -3255.125  # N20: -3775.5 -> -3255.125
$ws.Cells.Item(38, 8).Value = 15018  # H38: 21333.334 -> 15018
$ws.Cells.Item(38, 9).Value = 0  # I38: 29000 -> 0
$ws.Cells.Item(38, 10).Value = 15018  # J38: 17500 -> 15018
$ws.Cells.Item(38, 11).Value = 0  # K38: 29000 -> 0
$ws.Cells.Item(38, 12).ClearContents()  # L38: delete (was 17500)
$ws.Cells.Item(38, 13).Value = 15018  # M38: -28584 -> 15018
$ws.Cells.Item(38, 14).Value = -15850  # N38: -18332 -> -15850
$ws.Cells.Item(86, 8).Value = 2333.5557  # H86: 1840.8518 -> 2333.5557
$ws.Cells.Item(86, 9).Value = 2525.8462  # I86: 1947.75 -> 2525.8462
$ws.Cells.Item(86, 10).Value = 1833.6  # J86: 1535.4286 -> 1833.6
$ws.Cells.Item(86, 11).Value = 2525.8462  # K86: 1947.75 -> 2525.8462
$ws.Cells.Item(86, 12).Value = 1833.6  # L86: 1535.4286 -> 1833.6
$ws.Cells.Item(86, 13).Value = -1402.8462  # M86: -824.75 -> -1402.8462
$ws.Cells.Item(86, 14).Value = -4079.6  # N86: -3781.4286 -> -4079.6
$ws.Cells.Item(89, 8).Value = 2333.5557  # H89: 1840.8518 -> 2333.5557
$ws.Cells.Item(89, 9).Value = 2525.8462  # I89: 1947.75 -> 2525.8462
$ws.Cells.Item(89, 10).Value = 1833.6  # J89: 1535.4286 -> 1833.6
$ws.Cells.Item(89, 11).Value = 12629.231  # K89: 9738.75 -> 12629.231
$ws.Cells.Item(89, 12).Value = 9168  # L89: 7677.143 -> 9168
$ws.Cells.Item(89, 13).Value = -7013.231  # M89: -4122.75 -> -7013.231
$ws.Cells.Item(89, 14).Value = -20400  # N89: -18909.143 -> -20400

# ---- Sheet: CRP (index 4) ----
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(7, 8).Value = 75.347824  # H7: 75.39130400000001 -> 75.347824
$ws.Cells.Item(7, 9).Value = 53.75  # I7: 53.833332 -> 53.75
$ws.Cells.Item(7, 11).Value = 53.75  # K7: 53.833332 -> 53.75
$ws.Cells.Item(7, 13).Value = 59.25  # M7: 59.166668 -> 59.25
$ws.Cells.Item(31, 10).Value = 10000  # J31: 0 -> 10000
$ws.Cells.Item(31, 12).Value = 10000  # L31: 0 -> 10000
$ws.Cells.Item(31, 14).Value = -10590  # N31: None -> -10590
$ws.Cells.Item(34, 10).Value = 10000  # J34: 0 -> 10000
$ws.Cells.Item(34, 12).Value = 10000  # L34: 0 -> 10000
$ws.Cells.Item(34, 14).Value = -10404  # N34: None -> -10404
$ws.Cells.Item(132, 8).Value = 4656.771  # H132: 4755.6445 -> 4656.771
$ws.Cells.Item(132, 9).Value = 4315.268  # I132: 4382.1797 -> 4315.268
$ws.Cells.Item(132, 10).Value = 6657  # J132: 7183.1665 -> 6657
$ws.Cells.Item(132, 11).Value = 12945.804  # K132: 13146.5391 -> 12945.804
$ws.Cells.Item(132, 12).Value = 19971  # L132: 21549.4995 -> 19971
$ws.Cells.Item(132, 13).Value = -10415.804  # M132: -10616.5391 -> -10415.804
$ws.Cells.Item(132, 14).Value = -25031  # N132: -26609.4995 -> -25031

# ---- Sheet: CUL (index 5) ----
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(43, 8).Value = 12666  # H43: 20000 -> 12666
$ws.Cells.Item(43, 10).Value = 12666  # J43: 20000 -> 12666
$ws.Cells.Item(43, 12).Value = 37998  # L43: 60000 -> 37998
$ws.Cells.Item(43, 14).Value = -38226  # N43: -60228 -> -38226
$ws.Cells.Item(116, 8).Value = 6263.25  # H116: 5017.6665 -> 6263.25
$ws.Cells.Item(116, 9).Value = 6686  # I116: 5029 -> 6686
$ws.Cells.Item(116, 11).Value = 20058  # K116: 15087 -> 20058
$ws.Cells.Item(116, 13).Value = -16616  # M116: -11645 -> -16616

# ---- Sheet: GSM (index 6) ----
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(2, 8).Value = 4545525  # H2: 5000075 -> 4545525
$ws.Cells.Item(2, 9).Value = 5263208.5  # I2: 5882406.5 -> 5263208.5
$ws.Cells.Item(2, 11).Value = 5263208.5  # K2: 5882406.5 -> 5263208.5
$ws.Cells.Item(2, 13).Value = -5263095.5  # M2: -5882293.5 -> -5263095.5
$ws.Cells.Item(70, 8).Value = 15511.723  # H70: 15906.4 -> 15511.723
$ws.Cells.Item(70, 9).Value = 16443.867  # I70: 16952.346 -> 16443.867
$ws.Cells.Item(70, 11).Value = 16443.867  # K70: 16952.346 -> 16443.867
$ws.Cells.Item(70, 13).Value = -16173.867  # M70: -16682.346 -> -16173.867
$ws.Cells.Item(73, 8).Value = 15511.723  # H73: 15906.4 -> 15511.723
$ws.Cells.Item(73, 9).Value = 16443.867  # I73: 16952.346 -> 16443.867
$ws.Cells.Item(73, 11).Value = 16443.867  # K73: 16952.346 -> 16443.867
$ws.Cells.Item(73, 13).Value = -15507.867  # M73: -16016.346 -> -15507.867
$ws.Cells.Item(80, 8).Value = 2033.36  # H80: 1993.4814 -> 2033.36
$ws.Cells.Item(80, 9).Value = 1998.8667  # I80: 1939.5883 -> 1998.8667
$ws.Cells.Item(80, 11).Value = 1998.8667  # K80: 1939.5883 -> 1998.8667
$ws.Cells.Item(80, 13).Value = -1000.8667  # M80: -941.5882999999999 -> -1000.8667
$ws.Cells.Item(83, 8).Value = 2033.36  # H83: 1993.4814 -> 2033.36
$ws.Cells.Item(83, 9).Value = 1998.8667  # I83: 1939.5883 -> 1998.8667
$ws.Cells.Item(83, 11).Value = 9994.333500000001  # K83: 9697.941499999999 -> 9994.333500000001
$ws.Cells.Item(83, 13).Value = -5002.333500000001  # M83: -4705.941499999999 -> -5002.333500000001
$ws.Cells.Item(132, 8).Value = 13358.2  # H132: 14251.441 -> 13358.2
$ws.Cells.Item(132, 9).Value = 15576.655  # I132: 16740.535 -> 15576.655
$ws.Cells.Item(132, 11).Value = 46729.965  # K132: 50221.605 -> 46729.965
$ws.Cells.Item(132, 13).Value = -44199.965  # M132: -47691.605 -> -44199.965

# ---- Sheet: LTW (index 7) ----
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(16, 8).Value = 788.6842  # H16: 814.3889 -> 788.6842
$ws.Cells.Item(16, 9).Value = 776.94446  # I16: 803.4706 -> 776.94446
$ws.Cells.Item(16, 11).Value = 776.94446  # K16: 803.4706 -> 776.94446
$ws.Cells.Item(16, 13).Value = -606.94446  # M16: -633.4706 -> -606.94446
$ws.Cells.Item(22, 8).Value = 2854.92  # H22: 3047.652 -> 2854.92
$ws.Cells.Item(22, 10).Value = 3545.0588  # J22: 3932.6 -> 3545.0588
$ws.Cells.Item(22, 12).Value = 3545.0588  # L22: 3932.6 -> 3545.0588
$ws.Cells.Item(22, 14).Value = -4135.0588  # N22: -4522.6 -> -4135.0588
$ws.Cells.Item(27, 8).Value = 2854.92  # H27: 3047.652 -> 2854.92
$ws.Cells.Item(27, 10).Value = 3545.0588  # J27: 3932.6 -> 3545.0588
$ws.Cells.Item(27, 12).Value = 3545.0588  # L27: 3932.6 -> 3545.0588
$ws.Cells.Item(27, 14).Value = -3759.0588  # N27: -4146.6 -> -3759.0588
$ws.Cells.Item(46, 8).Value = 4578.048  # H46: 4424.5 -> 4578.048
$ws.Cells.Item(46, 10).Value = 4846.7896  # J46: 4664.45 -> 4846.7896
$ws.Cells.Item(46, 12).Value = 4846.7896  # L46: 4664.45 -> 4846.7896
$ws.Cells.Item(46, 14).Value = -5222.7896  # N46: -5040.45 -> -5222.7896
$ws.Cells.Item(48, 8).Value = 14900  # H48: 0 -> 14900
$ws.Cells.Item(48, 9).Value = 14900  # I48: 0 -> 14900
$ws.Cells.Item(48, 11).Value = 14900  # K48: 0 -> 14900
$ws.Cells.Item(48, 13).Value = -14239  # M48: None -> -14239
$ws.Cells.Item(55, 8).Value = 1069.3784  # H55: 1047.6316 -> 1069.3784
$ws.Cells.Item(55, 10).Value = 1200.1818  # J55: 1158.5652 -> 1200.1818
$ws.Cells.Item(55, 12).Value = 1200.1818  # L55: 1158.5652 -> 1200.1818
$ws.Cells.Item(55, 14).Value = -1546.1818  # N55: -1504.5652 -> -1546.1818
$ws.Cells.Item(130, 8).Value = 100000  # H130: 0 -> 100000
$ws.Cells.Item(130, 10).Value = 100000  # J130: 0 -> 100000
$ws.Cells.Item(130, 12).Value = 100000  # L130: 0 -> 100000
$ws.Cells.Item(130, 14).Value = -110040  # N130: None -> -110040
$ws.Cells.Item(132, 8).Value = 3228.7778  # H132: 2274.75 -> 3228.7778
$ws.Cells.Item(132, 9).Value = 2166.6667  # I132: 1346.909 -> 2166.6667
$ws.Cells.Item(132, 10).Value = 3759.8333  # J132: 4316 -> 3759.8333
$ws.Cells.Item(132, 11).Value = 6500.000100000001  # K132: 4040.727 -> 6500.000100000001
$ws.Cells.Item(132, 12).Value = 11279.4999  # L132: 12948 -> 11279.4999
$ws.Cells.Item(132, 13).Value = -3970.000100000001  # M132: -1510.727 -> -3970.000100000001
$ws.Cells.Item(132, 14).Value = -16339.4999  # N132: -18008 -> -16339.4999
$ws.Cells.Item(136, 8).Value = 3889.1052  # H136: 4399.6 -> 3889.1052
$ws.Cells.Item(136, 9).Value = 3681.5  # I136: 4250.4165 -> 3681.5
$ws.Cells.Item(136, 11).Value = 11044.5  # K136: 12751.2495 -> 11044.5
$ws.Cells.Item(136, 13).Value = -8494.5  # M136: -10201.2495 -> -8494.5

# ---- Sheet: WVR (index 8) ----
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(45, 8).Value = 12857.286  # H45: 14417.667 -> 12857.286
$ws.Cells.Item(45, 9).Value = 12623  # I45: 17187 -> 12623
$ws.Cells.Item(45, 11).Value = 12623  # K45: 17187 -> 12623
$ws.Cells.Item(45, 13).Value = -12132  # M45: -16696 -> -12132
$ws.Cells.Item(96, 8).Value = 1746.25  # H96: 1995 -> 1746.25
$ws.Cells.Item(96, 9).Value = 1746.25  # I96: 1995 -> 1746.25
$ws.Cells.Item(96, 11).Value = 1746.25  # K96: 1995 -> 1746.25
$ws.Cells.Item(96, 13).Value = -373.25  # M96: -622 -> -373.25
$ws.Cells.Item(100, 8).Value = 33334554  # H100: 33334624 -> 33334554
$ws.Cells.Item(100, 9).Value = 249.5  # I100: 424.5 -> 249.5
$ws.Cells.Item(100, 11).Value = 499  # K100: 849 -> 499
$ws.Cells.Item(100, 13).Value = 42  # M100: -308 -> 42
$ws.Cells.Item(107, 8).Value = 18183580  # H107: 18520312 -> 18183580
$ws.Cells.Item(107, 10).Value = 2865.9546  # J107: 2998.476 -> 2865.9546
$ws.Cells.Item(107, 12).Value = 8597.863799999999  # L107: 8995.428 -> 8597.863799999999
$ws.Cells.Item(107, 14).Value = -12437.8638  # N107: -12835.428 -> -12437.8638
$ws.Cells.Item(132, 8).Value = 2181.9167  # H132: 1998.1666 -> 2181.9167
$ws.Cells.Item(132, 9).Value = 1726.0667  # I132: 1603.3 -> 1726.0667
$ws.Cells.Item(132, 10).Value = 2941.6667  # J132: 2787.9 -> 2941.6667
$ws.Cells.Item(132, 11).Value = 5178.2001  # K132: 4809.9 -> 5178.2001
$ws.Cells.Item(132, 12).Value = 8825.000100000001  # L132: 8363.700000000001 -> 8825.000100000001
$ws.Cells.Item(132, 13).Value = -2648.2001  # M132: -2279.9 -> -2648.2001
$ws.Cells.Item(132, 14).Value = -13885.0001  # N132: -13423.7 -> -13885.0001
$ws.Cells.Item(136, 8).Value = 1407.2195  # H136: 1463.3158 -> 1407.2195
$ws.Cells.Item(136, 9).Value = 1316.1613  # I136: 1382.5358 -> 1316.1613
$ws.Cells.Item(136, 11).Value = 3948.4839  # K136: 4147.607400000001 -> 3948.4839
$ws.Cells.Item(136, 13).Value = -1398.4839  # M136: -1597.607400000001 -> -1398.4839
